# Apply the "training on harder mazes completed" update to the INIT-positions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (C1:BJ1): the maze/column index header used to run in descending
#     order (59 .. 0). After re-training it is stored in ascending order
#     (0 .. 59) instead, i.e. the same 60 values, reversed.
for ($col = 3; $col -le 62; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 3
}

# --- Labels that moved around in the shared-string table. Only the actual
#     displayed text matters here; the duplicate "Index 1 " entry collapsed
#     into "Index 1" and "Index 2" swapped with it.
$ws.Range("BK1").Value = "Index 1"
$ws.Range("A2").Value = "Index 2"

# --- Results table (BO:BR, rows 3-12): newly trained "Index 1" start-point
#     results. BQ/BR are formulas (=BO*0.05 / =BP*0.05) and recalc on their
#     own, so only the BO column itself needs new values; BP is unchanged.
$ws.Range("BO3").Value = 7
$ws.Range("BO4").Value = 21
$ws.Range("BO5").Value = 20
$ws.Range("BO6").Value = 40
$ws.Range("BO7").Value = 52
$ws.Range("BO8").Value = 49
$ws.Range("BO9").Value = 32
$ws.Range("BO10").Value = 20
$ws.Range("BO11").Value = 9
$ws.Range("BO12").Value = 6

# --- Active cell / selection moved to BK1 (top of the "Index" labels area).
$null = $ws.Activate()
$null = $ws.Range("BK1").Select()

$null = $wb.Save()
